# "major accuracy check update"
#
# s1cDNASample_hbrown_08.29.19.xlsx
#   - polyAIsolationProtocol (col G) kit code corrected to the "L" (large)
#     pack size for every sample row.
#   - s1Protocol (col H) catalog numbers un-flattened: every row now records
#     its own distinct lot/catalog number instead of repeating one value.
#   - roboticS1Prep (col I) switched from a hard-coded FALSE literal to a
#     live =FALSE() formula.
#   - column G widened to fit the longer kit name, and the active selection
#     moved from column I to column H to match the new data-entry focus.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 27

# 1. polyAIsolationProtocol (G2:G27): NEBNextPoly(A)E7490 -> NEBNextPoly(A)E7490L
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = "NEBNextPoly(A)E7490L"
}

# 2. s1Protocol (H2:H27): single shared "E7760" -> one distinct catalog
#    number per row (E7420, E7421, E7422, ...). The refreshed cells pick up
#    the sheet's base font instead of the legacy column formatting.
$ws.Cells.Item($firstRow, 8).Value = "E7420"
$ws.Cells.Item($firstRow, 8).Style = "Normal"
for ($r = 3; $r -le $lastRow; $r++) {
    $catalogNumber = 7421 + ($r - 3)
    $ws.Cells.Item($r, 8).Value = "E" + $catalogNumber
    $ws.Cells.Item($r, 8).Style = "Normal"
}

# 3. roboticS1Prep (I2:I27): static FALSE boolean -> live =FALSE() formula.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Formula = "=FALSE()"
}

# 4. Widen column G to fit the new kit name, leaving the other columns as-is.
$ws.Columns.Item(7).ColumnWidth = 26.5

# 5. Move the active selection from I2:I27 to H2:H27.
[void]$ws.Range("H2:H27").Select()
